$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 237.5
$ws.Range("I28").Value = 165.21053
$ws.Range("J28").Value = 433.7143
$ws.Range("K28").Value = 165.21053
$ws.Range("L28").Value = 433.7143
$ws.Range("M28").Value = 319.78947
$ws.Range("N28").Value = -1403.7143

$ws.Range("H92").Value = 33341526
$ws.Range("I92").Value = 50011790
$ws.Range("J92").Value = 992
$ws.Range("K92").Value = 50011790
$ws.Range("L92").Value = 992
$ws.Range("M92").Value = -50010542
$ws.Range("N92").Value = -3488

$ws.Range("H98").Value = 586.6
$ws.Range("I98").Value = 606.9474
$ws.Range("J98").Value = 200
$ws.Range("K98").Value = 606.9474
$ws.Range("L98").Value = 200
$ws.Range("M98").Value = 891.0526
$ws.Range("N98").Value = -3196

$ws.Range("H116").Value = 7780.05
$ws.Range("I116").Value = 14475.375
$ws.Range("J116").Value = 3316.5
$ws.Range("K116").Value = 14475.375
$ws.Range("L116").Value = 3316.5
$ws.Range("M116").Value = -11033.375
$ws.Range("N116").Value = -10200.5

$ws.Range("H122").Value = 586.6
$ws.Range("I122").Value = 606.9474
$ws.Range("J122").Value = 200
$ws.Range("K122").Value = 1820.8422
$ws.Range("L122").Value = 600
$ws.Range("M122").Value = 629.1578
$ws.Range("N122").Value = -5500

$ws.Range("H129").Value = 620.4194
$ws.Range("I129").Value = 488.82608
$ws.Range("J129").Value = 998.75
$ws.Range("K129").Value = 1466.47824
$ws.Range("L129").Value = 2996.25
$ws.Range("M129").Value = 3533.52176
$ws.Range("N129").Value = -12996.25

$ws.Range("H132").Value = 225881.84
$ws.Range("I132").Value = 3784.1853
$ws.Range("J132").Value = 559028.3
$ws.Range("K132").Value = 11352.5559
$ws.Range("L132").Value = 1677084.9
$ws.Range("M132").Value = -8822.555899999999
$ws.Range("N132").Value = -1682144.9

$ws.Range("H138").Value = 1315.96
$ws.Range("I138").Value = 660.3077
$ws.Range("J138").Value = 2026.25
$ws.Range("K138").Value = 1980.9231
$ws.Range("L138").Value = 6078.75
$ws.Range("M138").Value = 3159.0769
$ws.Range("N138").Value = -16358.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2574.3044
$ws.Range("I61").Value = 2700.4285
$ws.Range("J61").Value = 1250
$ws.Range("K61").Value = 2700.4285
$ws.Range("L61").Value = 1250
$ws.Range("M61").Value = -2488.4285
$ws.Range("N61").Value = -1674

$ws.Range("H74").Value = 18406.482
$ws.Range("I74").Value = 24980.523
$ws.Range("J74").Value = 1149.625
$ws.Range("K74").Value = 24980.523
$ws.Range("L74").Value = 1149.625
$ws.Range("M74").Value = -24106.523
$ws.Range("N74").Value = -2897.625

$ws.Range("H77").Value = 18406.482
$ws.Range("I77").Value = 24980.523
$ws.Range("J77").Value = 1149.625
$ws.Range("K77").Value = 124902.615
$ws.Range("L77").Value = 5748.125
$ws.Range("M77").Value = -120534.615
$ws.Range("N77").Value = -14484.125

$ws.Range("H132").Value = 3620594.2
$ws.Range("I132").Value = 5673579.5
$ws.Range("J132").Value = 777999.7
$ws.Range("K132").Value = 17020738.5
$ws.Range("L132").Value = 2333999.1
$ws.Range("M132").Value = -17018208.5
$ws.Range("N132").Value = -2339059.1

$ws.Range("H136").Value = 2574.3044
$ws.Range("I136").Value = 2700.4285
$ws.Range("J136").Value = 1250
$ws.Range("K136").Value = 8101.2855
$ws.Range("L136").Value = 3750
$ws.Range("M136").Value = -5551.2855
$ws.Range("N136").Value = -8850

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 209094.27
$ws.Range("I31").Value = 244850.31
$ws.Range("J31").Value = 21375
$ws.Range("K31").Value = 244850.31
$ws.Range("L31").Value = 21375
$ws.Range("M31").Value = -244555.31
$ws.Range("N31").Value = -21965

$ws.Range("H34").Value = 209094.27
$ws.Range("I34").Value = 244850.31
$ws.Range("J34").Value = 21375
$ws.Range("K34").Value = 244850.31
$ws.Range("L34").Value = 21375
$ws.Range("M34").Value = -244648.31
$ws.Range("N34").Value = -21779

$ws.Range("H35").Value = 15687.786
$ws.Range("I35").Value = 867.7273
$ws.Range("J35").Value = 70028
$ws.Range("K35").Value = 867.7273
$ws.Range("L35").Value = 70028
$ws.Range("M35").Value = -573.7273
$ws.Range("N35").Value = -70616

$ws.Range("H58").Value = 4745.92
$ws.Range("I58").Value = 1439.625
$ws.Range("J58").Value = 10623.777
$ws.Range("K58").Value = 1439.625
$ws.Range("L58").Value = 10623.777
$ws.Range("M58").Value = -1236.625
$ws.Range("N58").Value = -11029.777

$ws.Range("H94").Value = 2235.75
$ws.Range("I94").Value = 2620.111
$ws.Range("J94").Value = 1921.2727
$ws.Range("K94").Value = 2620.111
$ws.Range("L94").Value = 1921.2727
$ws.Range("M94").Value = -2169.111
$ws.Range("N94").Value = -2823.2727

$ws.Range("H132").Value = 1033.878
$ws.Range("I132").Value = 1048.3611
$ws.Range("J132").Value = 929.6
$ws.Range("K132").Value = 3145.0833
$ws.Range("L132").Value = 2788.8
$ws.Range("M132").Value = -615.0833000000002
$ws.Range("N132").Value = -7848.8

$ws.Range("H134").Value = 1509.12
$ws.Range("I134").Value = 1296
$ws.Range("J134").Value = 2057.1428
$ws.Range("K134").Value = 3888
$ws.Range("L134").Value = 6171.428400000001
$ws.Range("M134").Value = -1353
$ws.Range("N134").Value = -11241.4284

$ws.Range("H136").Value = 4745.92
$ws.Range("I136").Value = 1439.625
$ws.Range("J136").Value = 10623.777
$ws.Range("K136").Value = 4318.875
$ws.Range("L136").Value = 31871.331
$ws.Range("M136").Value = -1768.875
$ws.Range("N136").Value = -36971.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 166811500
$ws.Range("J121").Value = 214471760
$ws.Range("L121").Value = 643415280
$ws.Range("N121").Value = -643417900

$ws.Range("H131").Value = 28409912
$ws.Range("I131").Value = 437.5
$ws.Range("J131").Value = 34723130
$ws.Range("K131").Value = 1312.5
$ws.Range("L131").Value = 104169390
$ws.Range("M131").Value = 3727.5
$ws.Range("N131").Value = -104179470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 19497.582
$ws.Range("I132").Value = 1073.7715
$ws.Range("J132").Value = 51739.25
$ws.Range("K132").Value = 3221.3145
$ws.Range("L132").Value = 155217.75
$ws.Range("M132").Value = -691.3145000000004
$ws.Range("N132").Value = -160277.75

$ws.Range("H135").Value = 58166.668
$ws.Range("J135").Value = 58166.668
$ws.Range("L135").Value = 58166.668
$ws.Range("N135").Value = -68306.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 288515.34
$ws.Range("I132").Value = 73865.03999999999
$ws.Range("J132").Value = 717816
$ws.Range("K132").Value = 221595.12
$ws.Range("L132").Value = 2153448
$ws.Range("M132").Value = -219065.12
$ws.Range("N132").Value = -2158508

$ws.Range("H136").Value = 197461.36
$ws.Range("I136").Value = 271300.25
$ws.Range("J136").Value = 2315.7144
$ws.Range("K136").Value = 813900.75
$ws.Range("L136").Value = 6947.1432
$ws.Range("M136").Value = -811350.75
$ws.Range("N136").Value = -12047.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1619.7125
$ws.Range("I132").Value = 232.33333
$ws.Range("J132").Value = 5781.85
$ws.Range("K132").Value = 696.99999
$ws.Range("L132").Value = 17345.55
$ws.Range("M132").Value = 1833.00001
$ws.Range("N132").Value = -22405.55
